$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the stray "detail" cells in column B that no longer exist ---
# Row 4 currently holds "Patron" in A4 and "Gold, Silver, None" in B4.
# Row 9 currently holds "Payment " in A9 and "(Credit, cash, Check)" in B9.
$ws.Range("B4").Clear()
$ws.Range("B9").Clear()

# --- Step 2: insert the two new "Patron" detail rows right after row 4 ---
# (Insert() copies the formatting of the row above, which is exactly the
#  plain vertical-centred style used for rows 4-9 in the original sheet.)
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Range("B5").Clear()
$ws.Range("B6").Clear()

# --- Step 3: insert the two extra "Payment" detail rows ---
# After the inserts above, the old row 9 ("Payment ") is now row 11.
# Two more rows are needed immediately after it, before "OrderDetail".
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Range("B12").Clear()
$ws.Range("B13").Clear()

# --- Step 4: write the final text for every row in column A ---
$ws.Range("A1").Value = "Class Diagram List"
$ws.Range("A4").Value = "Patron "
$ws.Range("A5").Value = "Patron Gold"
$ws.Range("A6").Value = "Patron Silver"
$ws.Range("A7").Value = "FrontDesk"
$ws.Range("A8").Value = "Housekeeping"
$ws.Range("A9").Value = "RoomService"
$ws.Range("A10").Value = "Order"
$ws.Range("A11").Value = "Payment  Credit"
$ws.Range("A12").Value = "Payment Check"
$ws.Range("A13").Value = "Payment Cash"
$ws.Range("A14").Value = "OrderDetail"

# --- Step 5: restore the cursor/selection state recorded in the workbook ---
[void]$ws.Range("C13").Select()

$wb.Save()
